# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.746.86"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.633.16"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.80"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0636"
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.856.97"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.621.64"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.558"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.68"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "25.753.22"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.00"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.72"
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.123"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.87"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0493"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.901"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "1.127.19"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.546"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0156"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.57"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.78"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.805"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "1.766.68"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -3.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.04"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.416"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.60"
$ws.Range("E50").Value = "  -2.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.33"
$ws.Range("E51").Value = "  +2.70%  "
